$wb = $excel.ActiveWorkbook

# --- CashReceipt (sheet11): fix "ComapnyID" typo -> "CompanyID", update selection ---
$wsCashReceipt = $wb.Worksheets.Item("CashReceipt")
$wsCashReceipt.Range("F1").Value = "CompanyID"
$wsCashReceipt.Range("F1:F2").Select()

# --- CRATO_HomeCurrency (sheet12): remove the two rate columns (G:H), update selection ---
$wsHomeCurrency = $wb.Worksheets.Item("CRATO_HomeCurrency")
$wsHomeCurrency.Range("G1:H1").EntireColumn.Delete()
$wsHomeCurrency.Range("O15").Select()

# --- CRATO_ForeignCurrency (sheet13): remove the two rate columns (G:H), update selection ---
$wsForeignCurrency = $wb.Worksheets.Item("CRATO_ForeignCurrency")
$wsForeignCurrency.Range("G1:H1").EntireColumn.Delete()
$wsForeignCurrency.Range("O13").Select()

# --- CashReceipt_ForeignCurr (sheet14): add a new CompanyID column, update selection ---
# (selecting this last makes it the active tab, matching the target workbook view)
$wsCashReceiptForeign = $wb.Worksheets.Item("CashReceipt_ForeignCurr")
$wsCashReceiptForeign.Range("I1").Value = "CompanyID"
$wsCashReceiptForeign.Range("I2").Value = "aBb5f0000004JfX"
$wsCashReceiptForeign.Range("D5").Select()
